$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1293.2413
$ws.Range("I15").Value = 1293.2413
$ws.Range("K15").Value = 3879.7239
$ws.Range("M15").Value = -3710.7239

$ws.Range("H19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 1000
$ws.Range("N19").Value = -1350

$ws.Range("H33").Value = 343.17856
$ws.Range("I33").Value = 375.56
$ws.Range("J33").Value = 73.333336
$ws.Range("K33").Value = 375.56
$ws.Range("L33").Value = 73.333336
$ws.Range("M33").Value = -146.56
$ws.Range("N33").Value = -531.333336

$ws.Range("H80").Value = 8999
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 8999
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 26997
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -28993

$ws.Range("H83").Value = 8999
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 8999
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 80991
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -90975

$ws.Range("H112").Value = 2425.5833
$ws.Range("I112").Value = 1553.5
$ws.Range("K112").Value = 4660.5
$ws.Range("M112").Value = -3552.5

$ws.Range("H116").Value = 6729.8
$ws.Range("I116").Value = 7187.25
$ws.Range("J116").Value = 4900
$ws.Range("K116").Value = 7187.25
$ws.Range("L116").Value = 4900
$ws.Range("M116").Value = -3745.25
$ws.Range("N116").Value = -11784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2009.6666
$ws.Range("I2").Value = 2012.4286
$ws.Range("K2").Value = 2012.4286
$ws.Range("M2").Value = -1899.4286

$ws.Range("H32").Value = 10666.883
$ws.Range("I32").Value = 10500.774
$ws.Range("J32").Value = 12383.333
$ws.Range("K32").Value = 10500.774
$ws.Range("L32").Value = 12383.333
$ws.Range("M32").Value = -10213.774
$ws.Range("N32").Value = -12957.333

$ws.Range("H74").Value = 28987
$ws.Range("I74").Value = 28987
$ws.Range("K74").Value = 28987
$ws.Range("M74").Value = -28113

$ws.Range("H77").Value = 28987
$ws.Range("I77").Value = 28987
$ws.Range("K77").Value = 144935
$ws.Range("M77").Value = -140567

$ws.Range("H102").Value = 1435.7142
$ws.Range("I102").Value = 1435.7142
$ws.Range("K102").Value = 1435.7142
$ws.Range("M102").Value = 186.2858000000001

$ws.Range("H110").Value = 3733
$ws.Range("I110").Value = 1416.25
$ws.Range("K110").Value = 1416.25
$ws.Range("M110").Value = 628.75

$ws.Range("H116").Value = 2009.6666
$ws.Range("I116").Value = 2012.4286
$ws.Range("K116").Value = 2012.4286
$ws.Range("M116").Value = 281.5714

$ws.Range("H132").Value = 3181.3125
$ws.Range("I132").Value = 2761.8462
$ws.Range("K132").Value = 8285.5386
$ws.Range("M132").Value = -5755.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2009.6666
$ws.Range("I3").Value = 2012.4286
$ws.Range("K3").Value = 2012.4286
$ws.Range("M3").Value = -1898.4286

$ws.Range("H20").Value = 3796.4375
$ws.Range("I20").Value = 1604.4546
$ws.Range("K20").Value = 1604.4546
$ws.Range("M20").Value = -1357.4546

$ws.Range("H86").Value = 2756.4167
$ws.Range("I86").Value = 1759
$ws.Range("J86").Value = 3753.8333
$ws.Range("K86").Value = 1759
$ws.Range("L86").Value = 3753.8333
$ws.Range("M86").Value = -636
$ws.Range("N86").Value = -5999.8333

$ws.Range("H89").Value = 2756.4167
$ws.Range("I89").Value = 1759
$ws.Range("J89").Value = 3753.8333
$ws.Range("K89").Value = 8795
$ws.Range("L89").Value = 18769.1665
$ws.Range("M89").Value = -3179
$ws.Range("N89").Value = -30001.1665

$ws.Range("H94").Value = 4428.5
$ws.Range("I94").Value = 4428.5
$ws.Range("K94").Value = 4428.5
$ws.Range("M94").Value = -3977.5

$ws.Range("H105").Value = 2923.7144
$ws.Range("I105").Value = 2923.7144
$ws.Range("K105").Value = 2923.7144
$ws.Range("M105").Value = -1176.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3976.125
$ws.Range("I105").Value = 2968.1667
$ws.Range("K105").Value = 2968.1667
$ws.Range("M105").Value = -1221.1667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1283.6666
$ws.Range("I18").Value = 1481.3334
$ws.Range("J18").Value = 1217.7778
$ws.Range("K18").Value = 4444.0002
$ws.Range("L18").Value = 3653.3334
$ws.Range("M18").Value = -4275.0002
$ws.Range("N18").Value = -3991.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1006
$ws.Range("J80").Value = 1006
$ws.Range("L80").Value = 1006
$ws.Range("N80").Value = -3002

$ws.Range("H83").Value = 1006
$ws.Range("J83").Value = 1006
$ws.Range("L83").Value = 5030
$ws.Range("N83").Value = -15014

$ws.Range("H107").Value = 2290.818
$ws.Range("I107").Value = 2966.3333
$ws.Range("J107").Value = 1480.2
$ws.Range("K107").Value = 2966.3333
$ws.Range("L107").Value = 1480.2
$ws.Range("M107").Value = -1046.3333
$ws.Range("N107").Value = -5320.2

$ws.Range("H122").Value = 1571.5454
$ws.Range("I122").Value = 1328.7
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 3986.1
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -1536.1
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2483.1667
$ws.Range("I16").Value = 2483.1667
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2483.1667
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2313.1667
$ws.Range("N16").ClearContents()

$ws.Range("H68").Value = 2707.8948
$ws.Range("I68").Value = 2747.2222
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 2747.2222
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -1998.2222
$ws.Range("N68").Value = -3498

$ws.Range("H71").Value = 2707.8948
$ws.Range("I71").Value = 2747.2222
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 13736.111
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -9992.111000000001
$ws.Range("N71").Value = -17488

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H122").Value = 691.4
$ws.Range("I122").Value = 691.4
$ws.Range("K122").Value = 2074.2
$ws.Range("M122").Value = 375.8000000000002

$ws.Range("H126").Value = 1632.4286
$ws.Range("I126").Value = 1345.25
$ws.Range("K126").Value = 4035.75
$ws.Range("M126").Value = -1565.75
